$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 5: Assignee changed EM -> FL and the whole row re-styled to the
#     green "open issue" look used by the other rows (same pattern as row 8).
$ws.Range("A5").Interior.Color = $ws.Range("A3").Interior.Color
$ws.Range("A5").HorizontalAlignment = $ws.Range("A3").HorizontalAlignment

$ws.Range("B5").Interior.Color = $ws.Range("B3").Interior.Color
$ws.Range("B5").HorizontalAlignment = $ws.Range("B3").HorizontalAlignment

$ws.Range("C5").Interior.Color = $ws.Range("C3").Interior.Color
$ws.Range("C5").HorizontalAlignment = $ws.Range("C3").HorizontalAlignment
$ws.Range("C5").WrapText = $ws.Range("C3").WrapText

$ws.Range("D5").Interior.Color = $ws.Range("D3").Interior.Color
$ws.Range("D5").HorizontalAlignment = $ws.Range("D3").HorizontalAlignment

$ws.Range("E5").Interior.Color = $ws.Range("E3").Interior.Color
$ws.Range("E5").Value2 = "FL"

# --- New CDR rows (9-13) ------------------------------------------------
$ws.Range("A9").Value2 = "CDR"
$ws.Range("B9").Value2 = "Product Download"
$ws.Range("D9").Value2 = "Open"

$ws.Range("A10").Value2 = "CDR"
$ws.Range("B10").Value2 = "Browse URL"
$ws.Range("C10").Value2 = "What we should expect ?"
$ws.Range("D10").Value2 = "Open"

$ws.Range("A11").Value2 = "CDR"
$ws.Range("B11").Value2 = "Name/Description in DataSet Selection"
$ws.Range("D11").Value2 = "Open"

$ws.Range("A12").Value2 = "CDR"
$ws.Range("B12").Value2 = "Paging issue"
$ws.Range("C12").Value2 = "Automatic ? Manual ?"
$ws.Range("D12").Value2 = "Open"

$ws.Range("A13").Value2 = "CDR"
$ws.Range("B13").Value2 = "Granule representation in the Web Client"
$ws.Range("D13").Value2 = "Open"

# Filled in last, matching the shared-string order seen in the saved file.
$ws.Range("C9").Value2 = "Diffenriate between Direct Download and Download with Local DM"
$ws.Range("C9").HorizontalAlignment = $ws.Range("C2").HorizontalAlignment
$ws.Range("C9").WrapText = $ws.Range("C2").WrapText

$ws.Rows(9).RowHeight = 45

# --- Column B widened to fit the new longer titles ----------------------
$ws.Columns(2).ColumnWidth = 40.5

# --- Selection moved ------------------------------------------------------
[void]$ws.Range("B4").Select()
